$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N, shifting the existing
# N/O/P (Late / heading / Outstanding) columns right to O/P/Q.
$ws.Columns("N").EntireColumn.Insert()
$ws.Columns("N").ColumnWidth = 10.16666666666667

# Make "Repayment schedule" the active/selected sheet (was "NewLoanInput"),
# with cell K13 selected.
$ws.Activate() | Out-Null
$ws.Range("K13").Select() | Out-Null
